$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Row, NewPrice (D) or $null if unchanged, NewVolume (E)
$updates = @(
    @{ Row = 2; D = "30.355.49"; E = "  +0.62%  " }
    @{ Row = 3; D = "1.934.37"; E = "  +1.13%  " }
    @{ Row = 4; D = $null; E = "  +0.20%  " }
    @{ Row = 5; D = "251.93"; E = "  +2.45%  " }
    @{ Row = 6; D = "0.7208"; E = "  +2.78%  " }
    @{ Row = 7; D = $null; E = "  +0.01%  " }
    @{ Row = 8; D = "0.3285"; E = "  +1.95%  " }
    @{ Row = 9; D = "27.62"; E = "  +6.99%  " }
    @{ Row = 10; D = "0.07262"; E = "  +6.30%  " }
    @{ Row = 11; D = "0.8045"; E = "  +2.36%  " }
    @{ Row = 12; D = "0.08091"; E = "  +1.98%  " }
    @{ Row = 13; D = "1.934.11"; E = "  +1.10%  " }
    @{ Row = 14; D = $null; E = "  +1.26%  " }
    @{ Row = 15; D = "94.70"; E = "  +1.05%  " }
    @{ Row = 16; D = "15.03"; E = "  +5.08%  " }
    @{ Row = 17; D = "30.348.24"; E = "  +0.61%  " }
    @{ Row = 18; D = "253.73"; E = "  -1.97%  " }
    @{ Row = 19; D = "0.000008231"; E = "  +5.06%  " }
    @{ Row = 20; D = "5.806"; E = "  +0.63%  " }
    @{ Row = 21; D = "2.188.10"; E = "  +1.04%  " }
    @{ Row = 22; D = "1.001"; E = "  +0.02%  " }
    @{ Row = 23; D = "1.001"; E = "  +0.30%  " }
    @{ Row = 24; D = "6.944"; E = "  +2.32%  " }
    @{ Row = 25; D = "9.720"; E = "  +2.05%  " }
    @{ Row = 26; D = "166.24"; E = "  +4.80%  " }
    @{ Row = 27; D = "2.345"; E = "  +6.63%  " }
    @{ Row = 28; D = $null; E = "  +3.05%  " }
    @{ Row = 29; D = $null; E = "  -0.59%  " }
    @{ Row = 30; D = "1.359"; E = "  -0.44%  " }
    @{ Row = 31; D = $null; E = "  -0.01%  " }
    @{ Row = 32; D = "4.443"; E = "  +1.23%  " }
    @{ Row = 33; D = "4.193"; E = "  +0.82%  " }
    @{ Row = 34; D = $null; E = "  +4.34%  " }
    @{ Row = 35; D = "1.268"; E = "  +7.39%  " }
    @{ Row = 36; D = "0.7491"; E = "  +1.53%  " }
    @{ Row = 37; D = "2.765"; E = "  +1.43%  " }
    @{ Row = 38; D = "0.01967"; E = "  +2.89%  " }
    @{ Row = 39; D = "2.801"; E = "  +0.44%  " }
    @{ Row = 40; D = "79.23"; E = "  +0.38%  " }
    @{ Row = 41; D = "6.449"; E = "  -0.46%  " }
    @{ Row = 42; D = "0.4541"; E = "  +3.24%  " }
    @{ Row = 43; D = "2.031"; E = "  +1.70%  " }
    @{ Row = 44; D = "0.8433"; E = "  +1.76%  " }
    @{ Row = 45; D = $null; E = "  -0.06%  " }
    @{ Row = 46; D = "102.00"; E = "  +0.24%  " }
    @{ Row = 47; D = "9.798"; E = "  +2.41%  " }
    @{ Row = 48; D = "7.449"; E = "  +3.81%  " }
    @{ Row = 49; D = $null; E = "  +2.55%  " }
    @{ Row = 50; D = "0.4188"; E = "  +3.66%  " }
    @{ Row = 51; D = "0.06055"; E = "  +2.82%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Force text storage so numeric-looking strings (e.g. "30.355.49")
        # are not reinterpreted as numbers/dates by Excel.
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

Write-Output "Updated $($updates.Count) rows"
